$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New localization entries being "handed off":
#   eaf4563a-563d-4ade-8eb6-56a9aa18b524.md  (row 4 on every sheet)
#   f3feb4d7-776b-4f2d-adf2-2d4b76f883fb.md  (row 5 on every sheet)
# ---------------------------------------------------------------------

$commit = "9277b78b5b1e8193ed4f0cdb0d2ff33998788451"

$file1 = "eaf4563a-563d-4ade-8eb6-56a9aa18b524.md"
$file2 = "f3feb4d7-776b-4f2d-adf2-2d4b76f883fb.md"

$xlf1zh = "eaf4563a-563d-4ade-8eb6-56a9aa18b524.eebbe4f3fad91672387de7a5246c3c58c2f761ea.zh-cn.xlf"
$xlf2zh = "f3feb4d7-776b-4f2d-adf2-2d4b76f883fb.6be0d625f441fd7e039f4085aec6ff92d893db27.zh-cn.xlf"
$xlf1de = "eaf4563a-563d-4ade-8eb6-56a9aa18b524.eebbe4f3fad91672387de7a5246c3c58c2f761ea.de-de.xlf"
$xlf2de = "f3feb4d7-776b-4f2d-adf2-2d4b76f883fb.6be0d625f441fd7e039f4085aec6ff92d893db27.de-de.xlf"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(4,1).Value = $file1
$ov.Cells.Item(4,2).Value = "e2e\" + $file1
$ov.Cells.Item(4,3).Value = ".md"
$ov.Cells.Item(4,5).Value = "Ready for handoff"
$ov.Cells.Item(4,6).Value = "Ready for handoff"
$ov.Cells.Item(4,7).Value = "2017-02-22 06:41:28"

$ov.Cells.Item(5,1).Value = $file2
$ov.Cells.Item(5,2).Value = "e2e\" + $file2
$ov.Cells.Item(5,3).Value = ".md"
$ov.Cells.Item(5,5).Value = "Ready for handoff"
$ov.Cells.Item(5,6).Value = "Ready for handoff"
$ov.Cells.Item(5,7).Value = "2017-02-22 06:41:28"

$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/$commit/e2e/$file1", "", "", "e2e\" + $file1) | Out-Null
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/$commit/e2e/$file2", "", "", "e2e\" + $file2) | Out-Null

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G5"))

# ===========================================================================
# Sheets "zh-cn" and "de-de" share the same column layout:
#   A Source File Name      J Latest Target File
#   B File Extension        K Latest Handback File
#   C Status                L Latest Handback DateTime
#   D Source Path            M Latest Handback Name
#   E Priority               N Reference Tokens
#   F Content Duplicate      O To be localized
#   G Latest Handoff File    P Dependency From
#   H Latest Handoff Datetime Q Has metadata
#   I Lastest Handoff Name    R Error Detail
# ===========================================================================
function Fill-LangRow($ws, $row, $file, $xlfName, $datetime) {
    $ws.Cells.Item($row,1).Value = $file
    $ws.Cells.Item($row,2).Value = ".md"
    $ws.Cells.Item($row,3).Value = "Ready for handoff"
    $ws.Cells.Item($row,4).Value = "e2e"
    $ws.Cells.Item($row,5).Value = "ht"
    $ws.Cells.Item($row,6).Value = "'False"
    $ws.Cells.Item($row,7).Value = $xlfName
    $ws.Cells.Item($row,8).Value = $datetime
    $ws.Cells.Item($row,12).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item($row,15).Value = "'True"
    $ws.Cells.Item($row,17).Value = "'False"
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

Fill-LangRow $zh 4 $file1 $xlf1zh "2017-02-22 06:41:13"
Fill-LangRow $zh 5 $file2 $xlf2zh "2017-02-22 06:41:13"

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/$commit/e2e/$file1", "", "", $file1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/$commit/e2e/$file2", "", "", $file2) | Out-Null

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:R5"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

Fill-LangRow $de 4 $file1 $xlf1de "2017-02-22 06:41:28"
Fill-LangRow $de 5 $file2 $xlf2de "2017-02-22 06:41:28"

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/$commit/e2e/$file1", "", "", $file1) | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/$commit/e2e/$file2", "", "", $file2) | Out-Null

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:R5"))
